$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.553.34"
$ws.Range("E2").Value = "  +3.77%  "
$ws.Range("D3").Value = "3.493.69"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "590.56"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("D6").Value = "168.15"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.490.99"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  +7.55%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "0.127"
$ws.Range("E11").Value = "  +6.36%  "
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "4.097.52"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "28.01"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "66.574.06"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").Value = "0.0000177"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "3.489.83"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "390.11"
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("D22").Value = "7.89"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "72.75"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "0.531"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("E26").Value = "  +6.08%  "
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +7.90%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "6.32"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("D33").Value = "23.61"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").Value = "7.31"
$ws.Range("E34").Value = "  +4.55%  "
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("D36").Value = "162.84"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("D37").Value = "0.896"
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("D38").Value = "1.91"
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.80"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("D40").Value = "4.62"
$ws.Range("E40").Value = "  +6.64%  "
$ws.Range("D41").Value = "0.0738"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.30"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").Value = "2.782.75"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "26.51"
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("D45").Value = "42.71"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "2.53"
$ws.Range("E46").Value = "  +4.98%  "
$ws.Range("D47").Value = "0.0309"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "343.34"
$ws.Range("E48").Value = "  +4.88%  "
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").Value = "33.38"
$ws.Range("E50").Value = "  +11.50%  "
$ws.Range("D51").Value = "0.857"
$ws.Range("E51").Value = "  +5.87%  "
